$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B3:C17 values
$ws.Range("B3").Value = 2.2999999999999998
$ws.Range("C3").Value = 1.9

$ws.Range("B4").Value = 0.3
$ws.Range("C4").Value = 0.2

$ws.Range("B6").Value = 1.1000000000000001
$ws.Range("C6").Value = 0.4

$ws.Range("B7").Value = 2.2999999999999998
$ws.Range("C7").Value = 0.6

$ws.Range("B8").Value = 3.6
$ws.Range("C8").Value = 1.1000000000000001

$ws.Range("B9").Value = 5.9
$ws.Range("C9").Value = 1.8

$ws.Range("B10").Value = 7.5
$ws.Range("C10").Value = 2.2999999999999998

$ws.Range("B11").Value = 8.5
$ws.Range("C11").Value = 2.9

$ws.Range("B12").Value = 11.2
$ws.Range("C12").Value = 3.7

$ws.Range("B13").Value = 15.5
$ws.Range("C13").Value = 5.0999999999999996

$ws.Range("B14").Value = 22.1
$ws.Range("C14").Value = 7.7

$ws.Range("B15").Value = 32.6
$ws.Range("C15").Value = 11.4

$ws.Range("B16").Value = 40.9
$ws.Range("C16").Value = 16.600000000000001

$ws.Range("B17").Value = 89.5
$ws.Range("C17").Value = 67

# Update formulas in D3:D17 and E3:E17
$ws.Range("D3:D17").Formula = "=1-(B3/1000)"
$ws.Range("E3:E17").Formula = "=1-(C3/1000)"

# Update selection
$ws.Range("G15").Select()
